$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '[Feng%Pan%NULL%0,                            Lian%Yang%NULL%0,                            Yuncheng%Li%NULL%0,                            Bo%Liang%NULL%0,                            Lin%Li%NULL%0,                            Tianhe%Ye%NULL%0,                            Lingli%Li%NULL%0,                            Dehan%Liu%NULL%0,                            Shan%Gui%NULL%0,                            Yu%Hu%NULL%0,                            Chuansheng%Zheng%NULL%0]'
$ws.Range('E3').Value = '[Barnaby%D.P.%coreGivesNoEmail%0,                           Becker%L.B.%coreGivesNoEmail%0,                           Chelico%J.D.%coreGivesNoEmail%0,                           Cohen%S.L.%coreGivesNoEmail%0,                           Cookingham%J.%coreGivesNoEmail%0,                           Coppa%K.%coreGivesNoEmail%0,                           Crawford%J.M.%coreGivesNoEmail%0,                           Davidson%K.W.%coreGivesNoEmail%0,                           Diefenbach%M.A.%coreGivesNoEmail%0,                           Dominello%A.J.%coreGivesNoEmail%0,                           Duer-Hefele%J.%coreGivesNoEmail%0,                           Falzon%L.%coreGivesNoEmail%0,                           Gitlin%J.%coreGivesNoEmail%0,                           Hajizadeh%N.%coreGivesNoEmail%0,                           Harvin%T.G.%coreGivesNoEmail%0,                           Hirsch%J.S.%coreGivesNoEmail%0,                           Hirschwerk%D.A.%coreGivesNoEmail%0,                           Kim%E.J.%coreGivesNoEmail%0,                           Kozel%Z.M.%coreGivesNoEmail%0,                           Marrast%L.M.%coreGivesNoEmail%0,                           McGinn%T.%coreGivesNoEmail%0,                           Mogavero%J.N.%coreGivesNoEmail%0,                           Narasimhan%M.%coreGivesNoEmail%0,                           Osorio%G.A.%coreGivesNoEmail%0,                           Qiu%M.%coreGivesNoEmail%0,                           Richardson%S.%coreGivesNoEmail%0,                           Zanos%T.P.%coreGivesNoEmail%0]'
$ws.Range('E4').Value = '[Mario%Rivera-Izquierdo%NULL%0,                            María%del Carmen Valero-Ubierna%NULL%0,                            María%del Carmen Valero-Ubierna%NULL%0,                            Juan Luis%R-delAmo%NULL%0,                            Miguel Ángel%Fernández-García%NULL%0,                            Silvia%Martínez-Diz%NULL%0,                            Arezu%Tahery-Mahmoud%NULL%0,                            Marta%Rodríguez-Camacho%NULL%0,                            Ana Belén%Gámiz-Molina%NULL%0,                            Nicolás%Barba-Gyengo%NULL%0,                            Pablo%Gámez-Baeza%NULL%0,                            Celia%Cabrero-Rodríguez%NULL%0,                            Pedro Antonio%Guirado-Ruiz%NULL%0,                            Divina Tatiana%Martín-Romero%NULL%0,                            Antonio Jesús%Láinez-Ramos-Bossini%NULL%0,                            María Rosa%Sánchez-Pérez%NULL%0,                            José%Mancera-Romero%NULL%0,                            Miguel%García-Martín%NULL%0,                            Luis Miguel%Martín-delosReyes%NULL%0,                            Virginia%Martínez-Ruiz%NULL%0,                            Virginia%Martínez-Ruiz%NULL%0,                            Pablo%Lardelli-Claret%NULL%0,                            Eladio%Jiménez-Mejías%NULL%0,                            Muhammad%Adrish%NULL%0,                            Muhammad%Adrish%NULL%0,                            NULL%NULL%NULL%0,                            NULL%NULL%NULL%0]'
$ws.Range('E5').Value = '[Qiurong%Ruan%NULL%0,                            Kun%Yang%NULL%0,                            Kun%Yang%NULL%0,                            Wenxia%Wang%NULL%0,                            Wenxia%Wang%NULL%0,                            Lingyu%Jiang%NULL%0,                            Lingyu%Jiang%NULL%0,                            Jianxin%Song%songsingsjx@sina.com%0,                            Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range('E6').Value = '[Grace%Salacup%salacupg@einstein.edu%0,                            Kevin Bryan%Lo%NULL%0,                            Kevin Bryan%Lo%NULL%0,                            Fahad%Gul%NULL%0,                            Eric%Peterson%NULL%0,                            Robert%De Joy%NULL%0,                            Ruchika%Bhargav%NULL%0,                            Jerald%Pelayo%NULL%0,                            Jeri%Albano%NULL%0,                            Zurab%Azmaiparashvili%NULL%0,                            Sadia%Benzaquen%NULL%0,                            Gabriel%Patarroyo‐Aponte%NULL%0,                            Janani%Rangaswami%NULL%0]'
$ws.Range('E7').Value = '[Priyank%Shah%xref no email%0,       Jack%Owens%xref no email%0,       James%Franklin%xref no email%0,       Akshat%Mehta%xref no email%0,       William%Heymann%xref no email%0,       William%Sewell%xref no email%0,       Jennifer%Hill%xref no email%0,       Krista%Barfield%xref no email%0,       Rajkumar%Doshi%xref no email%0]'
$ws.Range('E8').Value = '[Yufeng%Shang%NULL%0,                            Tao%Liu%NULL%0,                            Yongchang%Wei%NULL%0,                            Jingfeng%Li%NULL%0,                            Liang%Shao%NULL%0,                            Minghui%Liu%NULL%0,                            Yongxi%Zhang%NULL%0,                            Zhigang%Zhao%NULL%0,                            Haibo%Xu%NULL%0,                            Zhiyong%Peng%NULL%0,                            Fuling%Zhou%zhoufuling@whu.edu.cn%0,                            Xinghuan%Wang%wangxinghuan@whu.edu.cn%0]'
$ws.Range('E9').Value = '[Shaobo%Shi%NULL%0,                            Mu%Qin%qinmuae@163.com%0,                            Yuli%Cai%NULL%0,                            Tao%Liu%NULL%0,                            Bo%Shen%NULL%0,                            Fan%Yang%NULL%0,                            Sheng%Cao%NULL%0,                            Xu%Liu%NULL%0,                            Xu%Liu%NULL%0,                            Yaozu%Xiang%NULL%0,                            Qinyan%Zhao%NULL%0,                            He%Huang%huanghe1977@whu.edu.cn%0,                            Bo%Yang%yybb112@whu.edu.cn%0,                            Congxin%Huang%NULL%0]'
$ws.Range('E10').Value = '[Rita de Cássia Menezes%Soares%NULL%0,                            Larissa Rodrigues%Mattos%NULL%0,                            Letícia Martins%Raposo%NULL%0]'
$ws.Range('E11').Value = '[Haiying%Sun%NULL%0,                            Ruoqi%Ning%NULL%0,                            Yu%Tao%NULL%0,                            Chong%Yu%NULL%0,                            Xiaoyan%Deng%NULL%0,                            Caili%Zhao%NULL%0,                            Silu%Meng%NULL%0,                            Fangxu%Tang%89650793@qq.com%0,                            Dong%Xu%89650793@qq.com%0,                            Dong%Xu%89650793@qq.com%0]'
$ws.Range('E12').Value = '[Kun%Wang%NULL%0,                            Peiyuan%Zuo%NULL%0,                            Yuwei%Liu%NULL%0,                            Meng%Zhang%NULL%0,                            Xiaofang%Zhao%NULL%0,                            Songpu%Xie%NULL%0,                            Hao%Zhang%NULL%0,                            Xinglin%Chen%NULL%0,                            Chengyun%Liu%chengyunliu@hust.edu.cn%0]'
$ws.Range('E13').Value = '[Bo%XU%NULL%0,                            Cun-yu%FAN%NULL%0,                            An-lu%WANG%NULL%0,                            Yi-long%ZOU%NULL%0,                            Yi-han%YU%NULL%0,                            Cong%HE%NULL%0,                            Wen-guang%XIA%NULL%0,                            Ji-xian%ZHANG%NULL%0,                            Qing%MIAO%NULL%0]'
$ws.Range('E14').Value = '[Xisheng%Yan%NULL%0,                            Fen%Li%NULL%0,                            Xiao%Wang%NULL%0,                            Jie%Yan%NULL%0,                            Fen%Zhu%NULL%0,                            Shifan%Tang%NULL%0,                            Yingzhong%Deng%NULL%0,                            Hua%Wang%NULL%0,                            Rui%Chen%NULL%0,                            Zhili%Yu%NULL%0,                            Yaping%Li%NULL%0,                            Jingzhou%Shang%NULL%0,                            Lingjun%Zeng%NULL%0,                            Jie%Zhao%NULL%0,                            Chaokun%Guan%NULL%0,                            Qiaomei%Liu%NULL%0,                            Haifeng%Chen%NULL%0,                            Wei%Gong%NULL%0,                            Xin%Huang%NULL%0,                            Yu‐Jiao%Zhang%NULL%0,                            Jianguang%Liu%NULL%0,                            Xiaoyan%Dong%NULL%0,                            Wen%Zheng%zhengwen12@mails.jlu.edu.cn%0,                            Shaoping%Nie%spnie@126.com%0,                            Dongsheng%Li%dongshengli196809@163.com%0]'
$ws.Range('E15').Value = '[Qing%Yang%NULL%0,                            Ying%Zhou%NULL%0,                            Xinrong%Wang%NULL%0,                            Shan%Gao%NULL%0,                            Yang%Xiao%NULL%0,                            Weiming%Zhang%NULL%0,                            Yi%Hu%huyizxyy@163.com%0,                            Yafei%Wang%wyf_527@163.com%0]'
$ws.Range('E16').Value = '[Xiaobo%Yang%NULL%0,                            Qingyu%Yang%NULL%0,                            Yaxin%Wang%NULL%0,                            Yongran%Wu%NULL%0,                            Jiqian%Xu%NULL%0,                            Yuan%Yu%NULL%0,                            You%Shang%NULL%0]'
$ws.Range('E17').Value = '[Wenjing%Ye%yewenjing@xinhuamed.com.cn%0,                            Guoxi%Chen%4799082@qq.com%0,                            Xiaopan%Li%xiaopanli0224@126.com%0,                            Xing%Lan%474137452@qq.com%0,                            Chen%Ji%c.ji.3@warwick.ac.uk%0,                            Min%Hou%min-hou710@hotmail.com%0,                            Di%Zhang%zhangdizhangdi1234@163.com%0,                            Guangwang%Zeng%785663609@qq.com%0,                            Yaling%Wang%953822402@qq.com%0,                            Cheng%Xu%xucheng112358@126.com%0,                            Weiwei%Lu%luweiwei100@126.com%0,                            Ruolin%Cui%crlqwerty@163.com%0,                            Yuyang%Cai%caiyuyang@sjtu.edu.cn%0,                            Hai%Huang%1220775601@qq.com%0,                            Ling%Yang%yangling01@xinhuamed.com.cn%0]'
$ws.Range('E18').Value = '[Caizheng%Yu%NULL%0,                            Qing%Lei%NULL%0,                            Wenkai%Li%NULL%0,                            Xiong%Wang%NULL%0,                            Wei%Liu%NULL%0,                            Xionglin%Fan%NULL%0,                            Wengang%Li%228907211@qq.com%0]'
$ws.Range('E19').Value = '[Fei%Zhou%NULL%0,                            Ting%Yu%NULL%0,                            Ronghui%Du%NULL%0,                            Guohui%Fan%NULL%0,                            Ying%Liu%NULL%0,                            Zhibo%Liu%NULL%0,                            Jie%Xiang%NULL%0,                            Yeming%Wang%NULL%0,                            Bin%Song%NULL%0,                            Xiaoying%Gu%NULL%0,                            Lulu%Guan%NULL%0,                            Yuan%Wei%NULL%0,                            Hui%Li%NULL%0,                            Xudong%Wu%NULL%0,                            Jiuyang%Xu%NULL%0,                            Shengjin%Tu%NULL%0,                            Yi%Zhang%NULL%0,                            Hua%Chen%NULL%0,                            Bin%Cao%NULL%0]'
$ws.Range('E20').Value = '[Moran%Amit%NULL%0,                            Alex%Sorkin%NULL%0,                            Jacob%Chen%NULL%0,                            Barak%Cohen%NULL%0,                            Barak%Cohen%NULL%0,                            Dana%Karol%NULL%0,                            Dana%Karol%NULL%0,                            Avishai M%Tsur%NULL%0,                            Shaul%Lev%NULL%0,                            Shaul%Lev%NULL%0,                            Tal%Rozenblat%NULL%0,                            Ayana%Dvir%NULL%0,                            Ayana%Dvir%NULL%0,                            Geva%Landau%NULL%0,                            Lidar%Fridrich%NULL%0,                            Lidar%Fridrich%NULL%0,                            Elon%Glassberg%NULL%0,                            Shani%Kesari%NULL%0,                            Sigal%Sviri%NULL%0,                            Ram%Gelman%NULL%0,                            Asaf%Miller%NULL%0,                            Danny%Epstein%NULL%0,                            Ronny%Ben-Avi%NULL%0,                            Ronny%Ben-Avi%NULL%0,                            Moshe%Matan%NULL%0,                            Daniel J.%Jakobson%NULL%0,                            Daniel J.%Jakobson%NULL%0,                            Tarif%Bader%NULL%0,                            Tarif%Bader%NULL%0,                            David%Dahan%NULL%0,                            David%Dahan%NULL%0,                            Daniel A.%King%NULL%0,                            Anat%Ben-Ari%NULL%0,                            Arie%Soroksky%NULL%0,                            Alon%Bar%NULL%0,                            Alon%Bar%NULL%0,                            Noam%Fink%NULL%0,                            Pierre%Singer%NULL%0,                            Avi%Benov%NULL%0]'
$ws.Range('E21').Value = '[Alexander%Muacevic%NULL%0,                            John R%Adler%NULL%0,                            Muhammad Sohaib%Asghar%NULL%0,                            Muhammad Sohaib%Asghar%NULL%0,                            Syed Jawad%Haider Kazmi%NULL%0,                            Noman%Ahmed Khan%NULL%0,                            Mohammed%Akram%NULL%0,                            Salman%Ahmed Khan%NULL%0,                            Uzma%Rasheed%NULL%0,                            Maira%Hassan%NULL%0,                            Gul Muhammad%Memon%NULL%0]'
$ws.Range('E22').Value = '[Pedro%Baqui%NULL%0,                            Ioana%Bica%NULL%0,                            Valerio%Marra%marra@cosmo-ufes.org%0,                            Ari%Ercole%NULL%0,                            Mihaela%van der Schaar%NULL%0]'
$ws.Range('E23').Value = '[Anesi%Adriano%coreGivesNoEmail%0,                           Bettinardi%Alessandra%coreGivesNoEmail%0,                           Bonetti%Graziella%coreGivesNoEmail%0,                           Borrelli%Gianluca%coreGivesNoEmail%0,                           Fiordalisi%Gianfranco%coreGivesNoEmail%0,                           Lippi%Giuseppe%coreGivesNoEmail%0,                           Manelli%Filippo%coreGivesNoEmail%0,                           Marino%Antonio%coreGivesNoEmail%0,                           Menolfi%Annamaria%coreGivesNoEmail%0,                           Patroni%Andrea%coreGivesNoEmail%0,                           Saggini%Sara%coreGivesNoEmail%0,                           Volpi%Roberta%coreGivesNoEmail%0]'
$ws.Range('E24').Value = '[Andrea%Borghesi%NULL%0,                            Angelo%Zigliani%NULL%0,                            Salvatore%Golemi%NULL%0,                            Nicola%Carapella%NULL%0,                            Patrizia%Maculotti%NULL%0,                            Davide%Farina%NULL%0,                            Roberto%Maroldi%NULL%0]'
$ws.Range('E25').Value = '[Alberto M.%Borobia%NULL%0,                            Antonio J.%Carcas%NULL%0,                            Antonio J.%Carcas%NULL%0,                            Francisco%Arnalich%NULL%0,                            Rodolfo%Álvarez-Sala%NULL%0,                            Rodolfo%Álvarez-Sala%NULL%0,                            Jaime%Monserrat-Villatoro%NULL%0,                            Manuel%Quintana%NULL%0,                            Juan Carlos%Figueira%NULL%0,                            Rosario M.%Torres Santos-Olmo%NULL%0,                            Julio%García-Rodríguez%NULL%0,                            Julio%García-Rodríguez%NULL%0,                            Alberto%Martín-Vega%NULL%0,                            Antonio%Buño%NULL%0,                            Elena%Ramírez%NULL%0,                            Gonzalo%Martínez-Alés%NULL%0,                            Gonzalo%Martínez-Alés%NULL%0,                            Nicolás%García-Arenzana%NULL%0,                            M. Concepción%Núñez%NULL%0,                            M. Concepción%Núñez%NULL%0,                            Milagros%Martí-de-Gracia%NULL%0,                            Francisco%Moreno Ramos%NULL%0,                            Francisco%Reinoso-Barbero%NULL%0,                            Alejandro%Martin-Quiros%NULL%0,                            Angélica%Rivera Núñez%NULL%0,                            Jesús%Mingorance%NULL%0,                            Carlos J.%Carpio Segura%NULL%0,                            Carlos J.%Carpio Segura%NULL%0,                            Daniel%Prieto Arribas%NULL%0,                            Esther%Rey Cuevas%NULL%0,                            Concepción%Prados Sánchez%NULL%0,                            Juan J.%Rios%NULL%0,                            Miguel A.%Hernán%NULL%0,                            Jesús%Frías%NULL%0,                            José R.%Arribas%NULL%0,                            NULL%NULL%NULL%0]'
$ws.Range('E26').Value = '[Simon E.%Brill%simon.brill@nhs.net%0,                            Hannah C.%Jarvis%NULL%0,                            Hannah C.%Jarvis%NULL%0,                            Ezgi%Ozcan%NULL%0,                            Thomas L. P.%Burns%NULL%0,                            Rabia A.%Warraich%NULL%0,                            Lisa J.%Amani%NULL%0,                            Amina%Jaffer%NULL%0,                            Stephanie%Paget%NULL%0,                            Anand%Sivaramakrishnan%NULL%0,                            Dean D.%Creer%NULL%0]'
$ws.Range('E27').Value = '[Jianlei%Cao%NULL%0,                            Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                            Wenlin%Cheng%NULL%0,                            Lei%Yu%NULL%0,                            Ya-Kun%Liu%NULL%0,                            Xiaoyong%Hu%NULL%0,                            Qiang%Liu%NULL%0]'
$ws.Range('E28').Value = '[B.%Carter%NULL%0,                            J.T.%Collins%NULL%0,                            F.%Barlow-Pay%NULL%0,                            F.%Rickard%NULL%0,                            E.%Bruce%NULL%0,                            A.%Verduri%NULL%0,                            T.J.%Quinn%NULL%0,                            E.%Mitchell%NULL%0,                            A.%Price%NULL%0,                            A.%Vilches-Moraga%NULL%0,                            M.J.%Stechman%NULL%0,                            R.%Short%NULL%0,                            A.%Einarsson%NULL%0,                            P.%Braude%NULL%0,                            S.%Moug%NULL%0,                            P.K.%Myint%NULL%0,                            J.%Hewitt%NULL%0,                            L.%Pearce%NULL%0,                            K.%McCarthy%NULL%0,                            C.%Davey%NULL%0,                            S.%Jones%NULL%0,                            K.%Lunstone%NULL%0,                            A.%Cavenagh%NULL%0,                            C.%Silver%NULL%0,                            T.%Telford%NULL%0,                            R.%Simmons%NULL%0,                            M.%Holloway%NULL%0,                            J.%Hesford%NULL%0,                            T.%El Jichi Mutasem%NULL%0,                            S.%Singh%NULL%0,                            D.%Paxton%NULL%0,                            W.%Harris%NULL%0,                            N.%Galbraith%NULL%0,                            E.%Bhatti%NULL%0,                            J.%Edwards%NULL%0,                            S.%Duffy%NULL%0,                            J.%Kelly%NULL%0,                            C.%Murphy%NULL%0,                            C.%Bisset%NULL%0,                            R.%Alexander%NULL%0,                            M.%Garcia%NULL%0,                            S.%Sangani%NULL%0,                            T.%Kneen%NULL%0,                            T.%Lee%NULL%0,                            A.%McGovern%NULL%0,                            G.%Guaraldi%NULL%0,                            E.%Clini%NULL%0]'
$ws.Range('E29').Value = '[Fuyang%Chen%NULL%0,                            Wenwu%Sun%NULL%0,                            Shengrong%Sun%NULL%0,                            Zhiyu%Li%lizhiyu@whu.edu.cn%0,                            Zhong%Wang%zhongwangchn@whu.edu.cn%0,                            Li%Yu%yuliwhzxyy@163.com%0,                            Li%Yu%yuliwhzxyy@163.com%0]'
$ws.Range('E30').Value = '[Ruchong%Chen%NULL%0,                            Ling%Sang%NULL%0,                            Mei%Jiang%NULL%0,                            Zhaowei%Yang%NULL%0,                            Nan%Jia%NULL%0,                            Wanyi%Fu%NULL%0,                            Jiaxing%Xie%NULL%0,                            Weijie%Guan%NULL%0,                            Wenhua%Liang%NULL%0,                            Zhengyi%Ni%NULL%0,                            Yu%Hu%NULL%0,                            Lei%Liu%NULL%0,                            Hong%Shan%NULL%0,                            Chunliang%Lei%NULL%0,                            Yixiang%Peng%NULL%0,                            Li%Wei%NULL%0,                            Yong%Liu%NULL%0,                            Yahua%Hu%NULL%0,                            Peng%Peng%NULL%0,                            Jianming%Wang%NULL%0,                            Jiyang%Liu%NULL%0,                            Zhong%Chen%NULL%0,                            Gang%Li%NULL%0,                            Zhijian%Zheng%NULL%0,                            Shaoqin%Qiu%NULL%0,                            Jie%Luo%NULL%0,                            Changjiang%Ye%NULL%0,                            Shaoyong%Zhu%NULL%0,                            Jinping%Zheng%NULL%0,                            Nuofu%Zhang%NULL%0,                            Yimin%Li%NULL%0,                            Jianxing%He%NULL%0,                            Jing%Li%NULL%0,                            Shiyue%Li%NULL%0,                            Nanshan%Zhong%NULL%0,                            NULL%NULL%NULL%0]'
$ws.Range('E31').Value = '[Tao%Chen%NULL%0,                            Di%Wu%NULL%0,                            Huilong%Chen%NULL%0,                            Weiming%Yan%NULL%0,                            Danlei%Yang%NULL%0,                            Guang%Chen%NULL%0,                            Ke%Ma%NULL%0,                            Dong%Xu%NULL%0,                            Haijing%Yu%NULL%0,                            Hongwu%Wang%NULL%0,                            Tao%Wang%NULL%0,                            Wei%Guo%NULL%0,                            Jia%Chen%NULL%0,                            Chen%Ding%NULL%0,                            Xiaoping%Zhang%NULL%0,                            Jiaquan%Huang%NULL%0,                            Meifang%Han%NULL%0,                            Shusheng%Li%NULL%0,                            Xiaoping%Luo%NULL%0,                            Jianping%Zhao%NULL%0,                            Qin%Ning%NULL%0]'
$ws.Range('E32').Value = '[Anying%Cheng%NULL%0,                            Liu%Hu%NULL%0,                            Yiru%Wang%NULL%0,                            Luyan%Huang%NULL%0,                            Lingxi%Zhao%NULL%0,                            Congcong%Zhang%NULL%0,                            Xiyue%Liu%NULL%0,                            Ranran%Xu%NULL%0,                            Feng%Liu%NULL%0,                            Jinping%Li%NULL%0,                            Dawei%Ye%NULL%0,                            Tao%Wang%NULL%0,                            Yongman%Lv%lvyongman@126.com%0,                            Qingquan%Liu%qqliutj@163.com%0]'
$ws.Range('E33').Value = '[Fabio%Ciceri%ciceri.fabio@hsr.it%0,                            Antonella%Castagna%NULL%0,                            Patrizia%Rovere-Querini%NULL%0,                            Francesco%De Cobelli%NULL%0,                            Annalisa%Ruggeri%NULL%0,                            Laura%Galli%NULL%0,                            Caterina%Conte%NULL%0,                            Rebecca%De Lorenzo%NULL%0,                            Andrea%Poli%NULL%0,                            Alberto%Ambrosio%NULL%0,                            Carlo%Signorelli%NULL%0,                            Eleonora%Bossi%NULL%0,                            Maria%Fazio%NULL%0,                            Cristina%Tresoldi%NULL%0,                            Sergio%Colombo%NULL%0,                            Giacomo%Monti%NULL%0,                            Efgeny%Fominskiy%NULL%0,                            Stefano%Franchini%NULL%0,                            Marzia%Spessot%NULL%0,                            Carlo%Martinenghi%NULL%0,                            Michele%Carlucci%NULL%0,                            Luigi%Beretta%NULL%0,                            Anna Maria%Scandroglio%NULL%0,                            Massimo%Clementi%NULL%0,                            Massimo%Locatelli%NULL%0,                            Moreno%Tresoldi%NULL%0,                            Paolo%Scarpellini%NULL%0,                            Gianvito%Martino%NULL%0,                            Emanuele%Bosi%NULL%0,                            Lorenzo%Dagna%NULL%0,                            Adriano%Lazzarin%NULL%0,                            Giovanni%Landoni%NULL%0,                            Alberto%Zangrillo%NULL%0]'
$ws.Range('E34').Value = '[Yan%Deng%NULL%0,                            Wei%Liu%NULL%0,                            Kui%Liu%NULL%0,                            Yuan-Yuan%Fang%NULL%0,                            Jin%Shang%NULL%0,                            Ling%Zhou%NULL%0,                            Ke%Wang%NULL%0,                            Fan%Leng%NULL%0,                            Shuang%Wei%NULL%0,                            Lei%Chen%NULL%0,                            Hui-Guo%Liu%NULL%0,                            Pei-Fang%Wei%NULL%0,                            Pei-Fang%Wei%NULL%0]'
$ws.Range('E35').Value = '[Rong-Hui%Du%NULL%0,                            Li-Rong%Liang%NULL%0,                            Cheng-Qing%Yang%NULL%0,                            Wen%Wang%NULL%0,                            Tan-Ze%Cao%NULL%0,                            Ming%Li%NULL%0,                            Guang-Yun%Guo%NULL%0,                            Juan%Du%NULL%0,                            Chun-Lan%Zheng%NULL%0,                            Qi%Zhu%NULL%0,                            Ming%Hu%NULL%0,                            Xu-Yan%Li%NULL%0,                            Peng%Peng%NULL%0,                            Huan-Zhong%Shi%NULL%0]'
$ws.Range('E36').Value = '[Shan%Gao%NULL%0,                            Fang%Jiang%NULL%0,                            Wei%Jin%NULL%0,                            Yuan%Shi%NULL%0,                            Leilei%Yang%NULL%0,                            Yanqiong%Xia%NULL%0,                            Linyan%Jia%NULL%0,                            Bo%Wang%NULL%0,                            Han%Lin%NULL%0,                            Yin%Cai%NULL%0,                            Zhengyuan%Xia%NULL%0,                            Jian%Peng%NULL%0]'
$ws.Range('E37').Value = '[Pedro David%Wendel Garcia%pedrodavid.wendelgarcia@usz.ch%0,                            Thierry%Fumeaux%thierry.fumeaux@ghol.ch%0,                            Philippe%Guerci%NULL%0,                            Dorothea Monika%Heuberger%dorotheamonika.heuberger@usz.ch%0,                            Jonathan%Montomoli%NULL%0,                            Ferran%Roche-Campo%NULL%0,                            Reto Andreas%Schuepbach%reto.schuepbach@usz.ch%0,                            Matthias Peter%Hilty%matthias.hilty@usz.ch%0,                            Mario%Alfaro Farias%NULL%0,                            Antoni%Margarit%NULL%0,                            Gerardo%Vizmanos-Lamotte%NULL%0,                            Thomas%Tschoellitsch%NULL%0,                            Jens%Meier%NULL%0,                            Francesco S.%Cardona%NULL%0,                            Josef%Skola%NULL%0,                            Lenka%Horakova%NULL%0,                            Hernan%Aguirre-Bermeo%NULL%0,                            Janina%Apolo%NULL%0,                            Emmanuel%Novy%NULL%0,                            Marie-Reine%Losser%NULL%0,                            Geoffrey%Jurkolow%NULL%0,                            Gauthier%Delahaye%NULL%0,                            Sascha%David%NULL%0,                            Tobias%Welte%NULL%0,                            Tobias%Wengenmayer%NULL%0,                            Dawid L.%Staudacher%NULL%0,                            Theodoros%Aslanidis%NULL%0,                            Barna%Babik%NULL%0,                            Anita%Korsos%NULL%0,                            Janos%Gal%NULL%0,                            Hermann%Csaba%NULL%0,                            Abele%Donati%NULL%0,                            Andrea%Carsetti%NULL%0,                            Fabrizio%Turrini%NULL%0,                            Maria Sole%Simonini%NULL%0,                            Roberto%Ceriani%NULL%0,                            Martina%Murrone%NULL%0,                            Emanuele%Rezoagli%NULL%0,                            Giovanni%Vitale%NULL%0,                            Alberto%Fogagnolo%NULL%0,                            Savino%Spadaro%NULL%0,                            Maddalena Alessandra%Wu%NULL%0,                            Chiara%Cogliati%NULL%0,                            Riccardo%Colombo%NULL%0,                            Emanuele%Catena%NULL%0,                            Francesca%Facondini%NULL%0,                            Antonella%Potalivo%NULL%0,                            Gianfilippo%Gangitano%NULL%0,                            Tiziana%Perin%NULL%0,                            Maria Grazia%Bocci%NULL%0,                            Massimo%Antonelli%NULL%0,                            Diederik%Gommers%NULL%0,                            Can%Ince%NULL%0,                            Eric%Mayor-Vázquez%NULL%0,                            Maria%Cruz%NULL%0,                            Martin%Delgado%NULL%0,                            Raquel Rodriguez%Garcia%NULL%0,                            Jorge%Gamez Zapata%NULL%0,                            Begoña%Zalba-Etayo%NULL%0,                            Herminia%Lozano-Gomez%NULL%0,                            Pedro%Castro%NULL%0,                            Adrian%Tellez%NULL%0,                            Adriana%Jacas%NULL%0,                            Guido%Muñoz%NULL%0,                            Rut%Andrea%NULL%0,                            Jose%Ortiz%NULL%0,                            Eduard%Quintana%NULL%0,                            Irene%Rovira%NULL%0,                            Enric%Reverter%NULL%0,                            Javier%Fernandez%NULL%0,                            Miquel%Ferrer%NULL%0,                            Joan R.%Badia%NULL%0,                            Arantxa%Lander Azcona%NULL%0,                            Jesus Escos%Orta%NULL%0,                            Philipp%Bühler%NULL%0,                            Silvio%Brugger%NULL%0,                            Daniel%Hofmaenner%NULL%0,                            Simone%Unseld%NULL%0,                            Frank%Ruschitzka%NULL%0,                            Mallory%Moret-Bochatay%NULL%0,                            Bernd%Yuen%NULL%0,                            Thomas%Hillermann%NULL%0,                            Hatem%Ksouri%NULL%0,                            Govind Oliver%Sridharan%NULL%0,                            Anette%Ristic%NULL%0,                            Michael%Sepulcri%NULL%0,                            Miodrag%Filipovic%NULL%0,                            Urs%Pietsch%NULL%0,                            Petra%Salomon%NULL%0,                            Iris%Drvaric%NULL%0,                            Peter%Schott%NULL%0,                            Severin%Urech%NULL%0,                            Adriana%Lambert%NULL%0,                            Lukas%Merki%NULL%0,                            Marcus%Laube%NULL%0,                            Frank%Hillgaertner%NULL%0,                            Marianne%Sieber%NULL%0,                            Alexander%Dullenkopf%NULL%0,                            Lina%Petersen%NULL%0,                            Serge%Grazioli%NULL%0,                            Peter C.%Rimensberger%NULL%0,                            Isabelle%Fleisch%NULL%0,                            Jerome%Lavanchy%NULL%0,                            Katharina%Marquardt%NULL%0,                            Karim%Shaikh%NULL%0,                            Hermann%Redecker%NULL%0,                            Michael%Stephan%NULL%0,                            Jan%Brem%NULL%0,                            Bjarte%Rogdo%NULL%0,                            Andre%Birkenmaier%NULL%0,                            Friederike%Meyer zu Bentrup%NULL%0,                            Patricia%Fodor%NULL%0,                            Pascal%Locher%NULL%0,                            Giovanni%Camen%NULL%0,                            Martin%Siegemund%NULL%0,                            Nuria%Zellweger%NULL%0,                            Marie-Madlen%Jeitziner%NULL%0,                            Beatrice%Jenni-Moser%NULL%0,                            Christian%Bürkle%NULL%0,                            Gian-Reto%Kleger%NULL%0,                            Marilene%Franchitti Laurent%NULL%0,                            Jean-Christophe%Laurent%NULL%0,                            Tomislav%Gaspert%NULL%0,                            Marija%Jovic%NULL%0,                            Michael%Studhalter%NULL%0,                            Christoph%Haberthuer%NULL%0,                            Roger F.%Lussman%NULL%0,                            Daniela%Selz%NULL%0,                            Didier%Naon%NULL%0,                            Romano%Mauri%NULL%0,                            Samuele%Ceruti%NULL%0,                            Julien%Marrel%NULL%0,                            Mirko%Brenni%NULL%0,                            Rolf%Ensner%NULL%0,                            Nadine%Gehring%NULL%0,                            Antje%Heise%NULL%0,                            Tobias%Huebner%NULL%0,                            Thomas A.%Neff%NULL%0,                            Sara%Cereghetti%NULL%0,                            Filippo%Boroli%NULL%0,                            Jerome%Pugin%NULL%0,                            Nandor%Marczin%NULL%0,                            Joyce%Wong%NULL%0]'
$ws.Range('E38').Value = '[Warren%Gavin%NULL%0,                            Elliott%Campbell%NULL%0,                            Syed-Adeel%Zaidi%NULL%0,                            Neha%Gavin%NULL%0,                            Lana%Dbeibo%NULL%0,                            Cole%Beeler%NULL%0,                            Kari%Kuebler%NULL%0,                            Ahmed%Abdel-Rahman%NULL%0,                            Mark%Luetkemeyer%NULL%0,                            Areeba%Kara%NULL%0]'
$ws.Range('E39').Value = '[Vijay%Gayam%vgayam@interfaithmedical.com%0,                            Muchi Ditah%Chobufo%NULL%0,                            Muchi Ditah%Chobufo%NULL%0,                            Mohamed A.%Merghani%NULL%0,                            Mohamed A.%Merghani%NULL%0,                            Shristi%Lamichhane%NULL%0,                            Pavani Reddy%Garlapati%NULL%0,                            Mark K.%Adler%NULL%0]'
$ws.Range('E40').Value = '[Cao%Y%coreGivesNoEmail%0,                           Imam%Z%coreGivesNoEmail%0,                           Lippi%G%coreGivesNoEmail%0,                           Oran%DP%coreGivesNoEmail%0,                           Shi%S%coreGivesNoEmail%0]'
$ws.Range('E41').Value = '[Hai%Hu%huhai@wchscu.cn%0,                            Ni%Yao%NULL%0,                            Ni%Yao%NULL%0,                            Yanru%Qiu%NULL%0,                            John H.%Burton%NULL%0,                            John H.%Burton%NULL%0]'
$ws.Range('E42').Value = '[Jiaofeng%Huang%NULL%0,                            Aiguo%Cheng%NULL%0,                            Rahul%Kumar%NULL%0,                            Yingying%Fang%NULL%0,                            Yingying%Fang%NULL%0,                            Gongping%Chen%NULL%0,                            Yueyong%Zhu%NULL%0,                            Su%Lin%sumer5129@fjmu.edu.cn%0,                            Su%Lin%sumer5129@fjmu.edu.cn%0]'
$ws.Range('E43').Value = '[Jong-moon%Hwang%NULL%0,                            Ju-Hyun%Kim%NULL%0,                            Jin-Sung%Park%NULL%0,                            Min Cheol%Chang%wheel633@ynu.ac.kr%0,                            Donghwi%Park%bdome@hanmail.net%0]'
$ws.Range('E44').Value = '[K.%Khalil%NULL%0,                            K.%Agbontaen%NULL%0,                            D.%McNally%NULL%0,                            A.%Love%NULL%0,                            S.%Mandalia%NULL%0,                            W.%Banya%NULL%0,                            E.%Starren%NULL%0,                            R.%Dhunnookchand%NULL%0,                            H.%Farne%NULL%0,                            R.%Morton%NULL%0,                            G.%Davies%NULL%0,                            O.%Orhan%NULL%0,                            D%Lai%NULL%0,                            M.%Nelson%NULL%0,                            P.L.%Shah%NULL%0,                            J.L.%Garner%Justin.garner@chelwest.nhs.uk%0]'
$ws.Range('E45').Value = '[Eyal%Klang%NULL%0,                            Gassan%Kassim%NULL%0,                            Shelly%Soffer%soffer.shelly@gmail.com%0,                            Robert%Freeman%NULL%0,                            Robert%Freeman%NULL%0,                            Matthew A.%Levin%NULL%0,                            Matthew A.%Levin%NULL%0,                            David L.%Reich%NULL%0]'
$ws.Range('E46').Value = '[Sandeep%Krishnan%NULL%0,                            Kinjal%Patel%NULL%0,                            Ronak%Desai%NULL%0,                            Anupam%Sule%NULL%0,                            Peter%Paik%NULL%0,                            Ashley%Miller%NULL%0,                            Alicia%Barclay%NULL%0,                            Adam%Cassella%NULL%0,                            Jon%Lucaj%NULL%0,                            Yvonne%Royster%NULL%0,                            Joffer%Hakim%NULL%0,                            Zulfiqar%Ahmed%NULL%0,                            Farhad%Ghoddoussi%NULL%0]'
$ws.Range('E48').Value = '[ Manisha%Bhutani%null%0,                  David M.%Foureau%null%0,                  Shebli%Atrash%null%0,                  Peter M.%Voorhees%null%0,                  Saad Z.%Usmani%null%0]'
$ws.Range('E49').Value = '[Wil%Lieberman-Cribbin%NULL%0,                            Joseph%Rapp%NULL%0,                            Naomi%Alpert%NULL%0,                            Stephanie%Tuminello%NULL%0,                            Emanuela%Taioli%NULL%0]'
$ws.Range('E50').Value = '[Q.%Liu%NULL%0,                            N. C.%Song%NULL%0,                            Z. K.%Zheng%NULL%0,                            J. S.%Li%NULL%0,                            S. K.%Li%NULL%0]'
$ws.Range('E51').Value = '[Hui%Long%NULL%0,                            Lan%Nie%NULL%0,                            Xiaochen%Xiang%NULL%0,                            Huan%Li%NULL%0,                            Xiaoli%Zhang%NULL%0,                            Xiaozhi%Fu%NULL%0,                            Hongwei%Ren%NULL%0,                            Wanxin%Liu%NULL%0,                            Qiang%Wang%wangqiang@wust.edu.cn%0,                            Qingming%Wu%wuhe9224@sina.com%0,                            Qingming%Wu%wuhe9224@sina.com%0]'
$ws.Range('E52').Value = '[Miao%Luo%xref no email%0,       Jing%Liu%xref no email%0,       Weiling%Jiang%xref no email%0,       Shuang%Yue%xref no email%0,       Huiguo%Liu%xref no email%0,       Shuang%Wei%xref no email%0]'
$ws.Range('E53').Value = '[Xiaomin%Luo%luoxiaomin04@163.com%0,                            Wei%Zhou%NULL%0,                            Xiaojie%Yan%NULL%0,                            Tangxi%Guo%NULL%0,                            Benchao%Wang%NULL%0,                            Hongxia%Xia%NULL%0,                            Lu%Ye%NULL%0,                            Jun%Xiong%NULL%0,                            Zongping%Jiang%NULL%0,                            Yu%Liu%NULL%0,                            Bicheng%Zhang%NULL%0,                            Weize%Yang%NULL%0]'
$ws.Range('E54').Value = '[Ying%Luo%NULL%0,                            Liyan%Mao%NULL%0,                            Xu%Yuan%NULL%0,                            Ying%Xue%NULL%0,                            Qun%Lin%NULL%0,                            Guoxing%Tang%NULL%0,                            Huijuan%Song%NULL%0,                            Feng%Wang%fengwang@tjh.tjmu.edu.cn%0,                            Ziyong%Sun%zysun@tjh.tjmu.edu.cn%0]'
$ws.Range('E55').Value = '[Ying%Luo%NULL%0,                            Ying%Xue%NULL%0,                            Liyan%Mao%NULL%0,                            Xu%Yuan%NULL%0,                            Qun%Lin%NULL%0,                            Guoxing%Tang%NULL%0,                            Huijuan%Song%NULL%0,                            Feng%Wang%NULL%0,                            Ziyong%Sun%NULL%0]'
$ws.Range('E56').Value = '[Chiara%Masetti%NULL%0,                            Elena%Generali%NULL%0,                            Francesca%Colapietro%NULL%0,                            Antonio%Voza%NULL%0,                            Maurizio%Cecconi%NULL%0,                            Antonio%Messina%NULL%0,                            Paolo%Omodei%NULL%0,                            Claudio%Angelini%NULL%0,                            Michele%Ciccarelli%NULL%0,                            Salvatore%Badalamenti%NULL%0,                            G. Walter%Canonica%NULL%0,                            Ana%Lleo%ana.lleo@humanitas.it%0,                            Alessio%Aghemo%NULL%0,                            Alessio%Aghemo%NULL%0,                            NULL%NULL%NULL%0]'
$ws.Range('E57').Value = '[Takahisa%Mikami%NULL%0,                            Hirotaka%Miyashita%NULL%0,                            Takayuki%Yamada%NULL%0,                            Matthew%Harrington%NULL%0,                            Daniel%Steinberg%NULL%0,                            Andrew%Dunn%NULL%0,                            Evan%Siau%Evan.Siau@mountsinai.org%0]'
$ws.Range('E58').Value = '[Alexis K.%Okoh%alexis.okoh@rwjbh.org%0,                            Christoph%Sossou%NULL%0,                            Christoph%Sossou%NULL%0,                            Neha S.%Dangayach%NULL%0,                            Sherin%Meledathu%NULL%0,                            Oluwakemi%Phillips%NULL%0,                            Corinne%Raczek%NULL%0,                            Michael%Patti%NULL%0,                            Nathan%Kang%NULL%0,                            Sameer A.%Hirji%NULL%0,                            Charles%Cathcart%NULL%0,                            Christian%Engell%NULL%0,                            Marc%Cohen%NULL%0,                            Sandhya%Nagarakanti%NULL%0,                            Eliahu%Bishburg%NULL%0,                            Harpreet S.%Grewal%NULL%0]'
